$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("H4").Value = 9
$ws.Range("I4").Value = 12.5
$ws.Range("J4").Value = 5.1
$ws.Range("L4").Value = 1.23
$ws.Range("N4").Value = 5.8
$ws.Range("O4").Value = 1.16
$ws.Range("P4").Value = 2.64
$ws.Range("Q4").Value = 1.44
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 2.2
$ws.Range("T4").Value = 1.82
$ws.Range("U4").Value = 1.98
$ws.Range("V4").Value = 1.08
$ws.Range("AA4").Value = 350
$ws.Range("AB4").Value = 13.5
$ws.Range("AF4").Value = 11.5
$ws.Range("AG4").Value = 12.5
$ws.Range("AI4").Value = 120
$ws.Range("AN4").Value = 5.1
$ws.Range("AO4").Value = 160

# Row 5
$ws.Range("G5").Value = 5.9
$ws.Range("H5").Value = 1.76
$ws.Range("J5").Value = 3.6
$ws.Range("K5").Value = 4

# Row 6
$ws.Range("F6").Value = 1.85
$ws.Range("G6").Value = 1.97
$ws.Range("H6").Value = 4.4
$ws.Range("I6").Value = 4.9
$ws.Range("J6").Value = 3.8
$ws.Range("K6").Value = 4.4
$ws.Range("P6").Value = 2.64
$ws.Range("Q6").Value = 1.51

# Row 7
$ws.Range("H7").Value = 1.81
$ws.Range("I7").Value = 2.02
$ws.Range("J7").Value = 3.6
$ws.Range("P7").Value = 1.84
$ws.Range("Q7").Value = 1.76

# Row 10
$ws.Range("Q10").Value = 1.81

# Row 11
$ws.Range("I11").Value = 1.92
$ws.Range("Q11").Value = 1.4

# Row 12
$ws.Range("F12").Value = 1.84
$ws.Range("G12").Value = 2.34
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 4.9
$ws.Range("J12").Value = 2.72
$ws.Range("K12").Value = 7.4
$ws.Range("Q12").Value = 1.61

# Row 13
$ws.Range("F13").Value = 3.65
$ws.Range("H13").Value = 1.85
$ws.Range("I13").Value = 1.99

# Row 14
$ws.Range("P14").Value = 2.56

# Row 15
$ws.Range("H15").Value = 1.68

# Row 16
$ws.Range("J16").Value = 3.25

# Row 17
$ws.Range("F17").Value = 1.78
$ws.Range("G17").Value = 1000
$ws.Range("H17").Value = 1.8
$ws.Range("I17").Value = 2.28
$ws.Range("J17").Value = 1.78
$ws.Range("K17").Value = 1000
$ws.Range("P17").Value = 1.39
$ws.Range("Q17").Value = 2.28

# Row 18
$ws.Range("F18").Value = 2.12
$ws.Range("G18").Value = 2.86
$ws.Range("H18").Value = 3.1
$ws.Range("I18").Value = 4.9
$ws.Range("J18").Value = 2.24
$ws.Range("K18").Value = 5
$ws.Range("P18").Value = 1.44
$ws.Range("Q18").Value = 2.22

# Row 19
$ws.Range("H19").Value = 1.96

# Row 20
$ws.Range("F20").Value = 1.8
$ws.Range("G20").Value = 1.81
$ws.Range("H20").Value = 5.5
$ws.Range("I20").Value = 5.6
$ws.Range("K20").Value = 3.75
$ws.Range("N20").Value = 3.1
$ws.Range("R20").Value = 1.26
$ws.Range("U20").Value = 1.8
$ws.Range("Z20").Value = 42
$ws.Range("AD20").Value = 22

# Row 21
$ws.Range("O21").Value = 1.45
$ws.Range("AF21").Value = 16
$ws.Range("AK21").Value = 34

# Row 23
$ws.Range("G23").Value = 16.5
$ws.Range("Q23").Value = 1.83

# Row 24
$ws.Range("F24").Value = 1.89
$ws.Range("G24").Value = 2.08
$ws.Range("P24").Value = 1.63
$ws.Range("Q24").Value = 2.02

# Row 26
$ws.Range("F26").Value = 3.45
$ws.Range("H26").Value = 2.42
$ws.Range("I26").Value = 2.64
$ws.Range("J26").Value = 2.98
$ws.Range("K26").Value = 3.3

# Row 27
$ws.Range("I27").Value = 2.78
$ws.Range("J27").Value = 2.86

# Row 30
$ws.Range("F30").Value = 1.78
$ws.Range("G30").Value = 2.28
$ws.Range("H30").Value = 3.4
$ws.Range("J30").Value = 3.35
